$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns per latest crypto snapshot.
# Cells whose new text looks like a plain number need NumberFormat forced to
# Text ("@") first, otherwise Excel auto-converts the string to a numeric value
# and the literal formatting (e.g. trailing zeros) would be lost.

$ws.Range("D2").Value = "43.937.22"
$ws.Range("E2").Value = "  +0.49%  "
$ws.Range("D3").Value = "2.262.03"
$ws.Range("E3").Value = "  -0.33%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.657"
$ws.Range("E5").Value = "  +3.92%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "233.19"
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "63.79"
$ws.Range("E7").Value = "  +0.30%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.449"
$ws.Range("E9").Value = "  +3.88%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0978"
$ws.Range("E10").Value = "  -4.73%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "57.73"
$ws.Range("E11").Value = "  +0.22%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "26.72"
$ws.Range("E12").Value = "  +2.65%  "
$ws.Range("E13").Value = "  +1.67%  "
$ws.Range("D14").Value = "2.598.10"
$ws.Range("E14").Value = "  -0.29%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.64"
$ws.Range("E15").Value = "  -0.44%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.15"
$ws.Range("E16").Value = "  +3.42%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.842"
$ws.Range("E17").Value = "  +2.16%  "
$ws.Range("D18").Value = "2.267.72"
$ws.Range("E18").Value = "  +0.36%  "
$ws.Range("D19").Value = "43.853.89"
$ws.Range("E19").Value = "  +0.60%  "
$ws.Range("D20").Value = "0.0₃0983"
$ws.Range("E20").Value = "  -1.21%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "73.87"
$ws.Range("E21").Value = "  +0.15%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.16"
$ws.Range("E22").Value = "  +0.94%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "249.79"
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("E24").Value = "  -0.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.69"
$ws.Range("E25").Value = "  +31.76%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.46"
$ws.Range("E26").Value = "  -2.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.31"
$ws.Range("E27").Value = "  -0.91%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.92"
$ws.Range("E28").Value = "  +0.24%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "174.10"
$ws.Range("E29").Value = "  +0.97%  "
$ws.Range("E30").Value = "  +4.30%  "
$ws.Range("E31").Value = "  -0.71%  "
$ws.Range("E32").Value = "  -0.22%  "
$ws.Range("E33").Value = "  +3.56%  "
$ws.Range("E34").Value = "  +5.15%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0685"
$ws.Range("E35").Value = "  -0.42%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.97"
$ws.Range("E36").Value = "  -2.18%  "
$ws.Range("E37").Value = "  -2.49%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.46"
$ws.Range("E38").Value = "  -5.20%  "
$ws.Range("E39").Value = "  -1.36%  "
$ws.Range("E40").Value = "  +3.00%  "
$ws.Range("E41").Value = "  +0.13%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.70"
$ws.Range("E42").Value = "  +3.27%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.49"
$ws.Range("E43").Value = "  +0.73%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "17.26"
$ws.Range("E44").Value = "  -0.48%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "98.76"
$ws.Range("E45").Value = "  +1.05%  "
$ws.Range("E46").Value = "  -0.83%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.20"
$ws.Range("E47").Value = "  -1.01%  "
$ws.Range("E48").Value = "  +5.27%  "
$ws.Range("D49").Value = "1.457.31"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "10.01"
$ws.Range("E50").Value = "  -3.78%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.32"
$ws.Range("E51").Value = "  -1.19%  "
